# Insert a new data row at row 326 (shifting the existing rows 326-344 down
# to 327-345) and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 326; this pushes the old row 326
# (and everything below it) down by one row, and the new row inherits
# the number format (date format on column D) from the row above, just
# like a normal Excel "Insert Row" operation.
$ws.Rows.Item(326).Insert()

# Populate the new row with the new record's data.
$ws.Range("A326").Value = 6
$ws.Range("B326").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C326").Value = "Metropolitana"
$ws.Range("D326").Value = 45147
$ws.Range("E326").Value = 13
$ws.Range("F326").Value = 100112001
$ws.Range("G326").Value = "Berenjena"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Primera"
$ws.Range("J326").Value = 520
$ws.Range("K326").Value = 6000
$ws.Range("L326").Value = 7000
$ws.Range("M326").Value = 6385
$ws.Range("N326").Value = "`$/caja 50 unidades"
$ws.Range("O326").Value = "Región de Arica y Parinacota"
$ws.Range("P326").Value = 128
$ws.Range("Q326").Value = 50
$ws.Range("R326").Value = "Hortaliza"
